# luban excel export fixed
# Adds a new "BuffConfig" row (row 17) to the __tables__ worksheet,
# mirroring the existing SkillConfig row (row 16) pattern, and moves the
# active selection to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row describing the BuffConfig table export.
$ws.Range("B17").Value = "BuffConfigCategory"
$ws.Range("C17").Value = "BuffConfig"
$ws.Range("D17").Value = $true
$ws.Range("E17").Value = "BuffConfig.xlsx"

# Move/scroll selection to the newly added row, like a user would after
# typing the new entry.
$ws.Range("B17").Select()
